# feat: add 2022-Q3 data
#
# Insert a new worksheet "2022-Q3" right after "总计" (pushing 2022-Q2,
# 2022-Q1 and 2021-Q4 one slot to the right), populate it with the new
# quarter's fund-holding detail, and update the "总计" summary sheet with
# the new quarter's roll-up row.

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right after "总计"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"
$q3.Outline.SummaryRow = 1
$q3.Outline.SummaryColumn = 1

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q3.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
}
# header style (bold / centered / bordered) copied from the "总计" header
$totalSheet.Cells.Item(1, 2).Copy()
$q3.Range($q3.Cells.Item(1,2), $q3.Cells.Item(1,8)).PasteSpecial(-4122)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# code, name, scale, stock-position, ratio, market-value, rank
$rows = @(
    @("003853","金鹰信息产业股票A","8.78","88.43","4.23","0.3714",7),
    @("013233","华夏中证500指数智选增强A","21.06","93.71","1.09","0.2296",8),
    @("007994","华夏中证500指数增强A","19.56","93.37","1.11","0.2171",7),
    @("001637","嘉实量化精选股票","14.37","90.23","1.45","0.2084",10),
    @("005885","金鹰信息产业股票C","4.29","88.43","4.23","0.1815",7),
    @("016950","鹏华睿投灵活配置混合C","4.12","83.97","3.12","0.1285",3),
    @("013641","博道成长智航股票A","10.27","90.29","0.85","0.0873",4),
    @("013642","博道成长智航股票C","7.24","90.29","0.85","0.0615",4),
    @("007995","华夏中证500指数增强C","5.27","93.37","1.11","0.0585",7),
    @("013250","红土创新智能制造混合","1.31","93.66","3.71","0.0486",10),
    @("013234","华夏中证500指数智选增强C","3.92","93.71","1.09","0.0427",8),
    @("006441","中信建投中证500指数增强C","2.25","93.60","0.87","0.0196",9),
    @("006440","中信建投中证500指数增强A","2.14","93.60","0.87","0.0186",9),
    @("011590","九泰天利量化股票C","0.50","83.77","2.43","0.0122",2),
    @("004481","华宝第三产业灵活配置混合A","0.60","86.50","1.70","0.0102",7),
    @("005260","银华稳健增利灵活配置混合A","0.28","91.67","0.74","0.0021",9),
    @("011589","九泰天利量化股票A","0.07","83.77","2.43","0.0017",2),
    @("008838","德邦量化对冲策略灵活配置混合A","0.16","69.90","0.92","0.0015",8),
    @("012798","华宝第三产业灵活配置混合C","0.02","86.50","1.70","0.0003",7),
    @("008839","德邦量化对冲策略灵活配置混合C","0.03","69.90","0.92","0.0003",8),
    @("005261","银华稳健增利灵活配置混合C","0.02","91.67","0.74","0.0001",9),
    @("005434","鹏华睿投灵活配置混合A","0.00","83.97","3.12",$null,3)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    # index column -- numeric, styled like the "总计" index column
    $totalSheet.Cells.Item(2, 1).Copy()
    $q3.Cells.Item($r, 1).PasteSpecial(-4122)
    $q3.Cells.Item($r, 1).Value = $i

    # fund code -- always text so leading zeros survive
    $q3.Cells.Item($r, 2).Value = "'" + $row[0]
    $q3.Cells.Item($r, 2).Style = "Normal"

    # fund name -- plain text
    $q3.Cells.Item($r, 3).Value = "'" + $row[1]
    $q3.Cells.Item($r, 3).Style = "Normal"

    # scale / stock-position / ratio -- numeric-looking text, kept as text
    $q3.Cells.Item($r, 4).Value = "'" + $row[2]
    $q3.Cells.Item($r, 4).Style = "Normal"
    $q3.Cells.Item($r, 5).Value = "'" + $row[3]
    $q3.Cells.Item($r, 5).Style = "Normal"
    $q3.Cells.Item($r, 6).Value = "'" + $row[4]
    $q3.Cells.Item($r, 6).Style = "Normal"

    # market value -- text, except the one true-zero row which is numeric
    if ($row[5] -eq $null) {
        $q3.Cells.Item($r, 7).Value = 0
    } else {
        $q3.Cells.Item($r, 7).Value = "'" + $row[5]
        $q3.Cells.Item($r, 7).Style = "Normal"
    }

    # rank -- numeric
    $q3.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# 2. Update "总计": shift the existing roll-up rows down one and add the
#    new 2022-Q3 roll-up as the new first data row.
# ---------------------------------------------------------------------

# shift 2021-Q4 (row4 -> row5)
$totalSheet.Cells.Item(4, 1).Copy()
$totalSheet.Cells.Item(5, 1).PasteSpecial(-4122)
$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(5, 2).Value = "'2021-Q4"
$totalSheet.Cells.Item(5, 2).Style = "Normal"
$totalSheet.Cells.Item(5, 3).Value = 7
$totalSheet.Cells.Item(5, 4).Value = 0.54

# shift 2022-Q1 (row3 -> row4)
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(4, 2).Value = "'2022-Q1"
$totalSheet.Cells.Item(4, 2).Style = "Normal"
$totalSheet.Cells.Item(4, 3).Value = 4
$totalSheet.Cells.Item(4, 4).Value = 0.03

# shift 2022-Q2 (row2 -> row3)
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 2).Value = "'2022-Q2"
$totalSheet.Cells.Item(3, 2).Style = "Normal"
$totalSheet.Cells.Item(3, 3).Value = 35
$totalSheet.Cells.Item(3, 4).Value = 4.2

# new 2022-Q3 (row2)
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "'2022-Q3"
$totalSheet.Cells.Item(2, 2).Style = "Normal"
$totalSheet.Cells.Item(2, 3).Value = 22
$totalSheet.Cells.Item(2, 4).Value = 1.7
